# Femacal de La Calera - Zanahoria: weekly update.
# A new price-report row (14 columns worth of data, dated 2021-11-29 /
# serial 44529) is inserted above the existing row 154, pushing all the
# following records (old rows 154-244) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 154; everything currently at/after row 154
# shifts down to make room (old row 154 becomes row 155, ..., old row 244
# becomes row 245).
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row 154 with the new weekly record.
$ws.Cells.Item(154, 1).Value  = 3
$ws.Cells.Item(154, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(154, 3).Value  = "Coquimbo"
$ws.Cells.Item(154, 4).Value  = 44529
$ws.Cells.Item(154, 5).Value  = 5
$ws.Cells.Item(154, 6).Value  = 100114013
$ws.Cells.Item(154, 7).Value  = "Zanahoria"
$ws.Cells.Item(154, 8).Value  = "Sin especificar"
$ws.Cells.Item(154, 9).Value  = "Primera"
$ws.Cells.Item(154, 10).Value = 340
$ws.Cells.Item(154, 11).Value = 6000
$ws.Cells.Item(154, 12).Value = 6500
$ws.Cells.Item(154, 13).Value = 6265
$ws.Cells.Item(154, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(154, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(154, 16).Value = 313
$ws.Cells.Item(154, 17).Value = 20
$ws.Cells.Item(154, 18).Value = "Hortaliza"
